$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 25.97000000000062
$ws.Range("H2").Value = 0.00000003462542164101023
$ws.Range("I2").Value = 0.00000003462542164101023
$ws.Range("L2").Value = 60.66627607951715
$ws.Range("M2").Value = "[36.77702247559904, 84.55552968343525]"
$ws.Range("N2").Value = 0.000006271387400502704
$ws.Range("O2").Value = 0.000006271387400502704
$ws.Range("P2").Value = 2.182447749340658
$ws.Range("Q2").Value = "[1.792500312859043, 2.5723951858222724]"
$ws.Range("R2").Value = 0.0000000000000106581410364015
$ws.Range("S2").Value = 0.0000000000000106581410364015
$ws.Range("T2").Value = 60.39289084692335
$ws.Range("U2").Value = "[47.912051322166306, 72.8737303716804]"
$ws.Range("V2").Value = 0.000000000001157296480869263
$ws.Range("W2").Value = 0.000000000001157296480869263
$ws.Range("X2").Value = 16.94938938938979
$ws.Range("Y2").Value = 15.33763763763801
$ws.Range("Z2").Value = 18.56114114114158

# Row 3
$ws.Range("F3").Value = 25.97000000000062
$ws.Range("H3").Value = 0.0001162856580482208
$ws.Range("I3").Value = 0.0001162856580482208
$ws.Range("L3").Value = 51.08778163160643
$ws.Range("M3").Value = "[23.120853840658086, 79.05470942255478]"
$ws.Range("N3").Value = 0.0006227802674536598
$ws.Range("O3").Value = 0.0006227802674536598
$ws.Range("P3").Value = 1.289342330302117
$ws.Range("Q3").Value = "[0.6729738016698859, 1.9057108589343477]"
$ws.Range("R3").Value = 0.0001193954919762596
$ws.Range("S3").Value = 0.0001193954919762596
$ws.Range("T3").Value = 60.6148913010818
$ws.Range("U3").Value = "[44.89977636708264, 76.33000623508097]"
$ws.Range("V3").Value = 0.0000000007501173016066787
$ws.Range("W3").Value = 0.0000000007501173016066787
$ws.Range("X3").Value = 20.64082082082132
$ws.Range("Y3").Value = 18.09321321321365
$ws.Range("Z3").Value = 23.18842842842898

# Row 4
$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 25.97000000000062
$ws.Range("H4").Value = 0.00005092953192187277
$ws.Range("I4").Value = 0.00005092953192187277
$ws.Range("L4").Value = 67.90801604408148
$ws.Range("M4").Value = "[31.35186700813245, 104.46416508003051]"
$ws.Range("N4").Value = 0.000516066722537456
$ws.Range("O4").Value = 0.000516066722537456
$ws.Range("P4").Value = 0.735868549489501
$ws.Range("Q4").Value = "[0.2201316173686534, 1.2516054816103486]"
$ws.Range("R4").Value = 0.006170650697357782
$ws.Range("S4").Value = 0.006170650697357782
$ws.Range("T4").Value = 70.00408831373407
$ws.Range("U4").Value = "[51.09280786407818, 88.91536876338996]"
$ws.Range("V4").Value = 0.00000000216102713501698
$ws.Range("W4").Value = 0.00000000216102713501698
$ws.Range("X4").Value = 22.92846846846902
$ws.Range("Y4").Value = 20.79679679679729
$ws.Range("Z4").Value = 25.06014014014074

# Row 5
$ws.Range("F5").Value = 25.97000000000062
$ws.Range("H5").Value = 0.00008442234152761241
$ws.Range("I5").Value = 0.00008442234152761241
$ws.Range("L5").Value = 58.34801399331155
$ws.Range("M5").Value = "[29.238736204734167, 87.45729178188893]"
$ws.Range("N5").Value = 0.0002077543971854556
$ws.Range("O5").Value = 0.0002077543971854556
$ws.Range("P5").Value = 0.3585000625718084
$ws.Range("Q5").Value = "[-0.22013161736865428, 0.9371317425122712]"
$ws.Range("R5").Value = 0.2185337462402988
$ws.Range("S5").Value = 0.2185337462402988
$ws.Range("T5").Value = 68.84336299860581
$ws.Range("U5").Value = "[52.331235342062996, 85.35549065514863]"
$ws.Range("V5").Value = 0.00000000009168532599801438
$ws.Range("W5").Value = 0.00000000009168532599801438
$ws.Range("X5").Value = 24.48822822822881
$ws.Range("Y5").Value = 22.09659659659712
$ws.Range("Z5").Value = 26.8798598598605

# Row 6
$ws.Range("B6").Value = 0
$ws.Range("F6").Value = 25.97000000000062
$ws.Range("H6").Value = 0.001219809814260397
$ws.Range("I6").Value = 0.001219809814260397
$ws.Range("L6").Value = 43.38461606112249
$ws.Range("M6").Value = "[15.713139427331598, 71.05609269491337]"
$ws.Range("N6").Value = 0.002837384147732269
$ws.Range("O6").Value = 0.002837384147732269
$ws.Range("P6").Value = -0.2390000417145384
$ws.Range("Q6").Value = "[-0.9685791164220783, 0.49057903299300154]"
$ws.Range("R6").Value = 0.5127499295517957
$ws.Range("S6").Value = 0.5127499295517957
$ws.Range("T6").Value = 67.76964284941288
$ws.Range("U6").Value = "[52.42061741338472, 83.11866828544103]"
$ws.Range("V6").Value = 0.00000000001796940374276801
$ws.Range("W6").Value = 0.00000000001796940374276801
$ws.Range("X6").Value = 0.9878478478478705
$ws.Range("Y6").Value = -2.027687687687739
$ws.Range("Z6").Value = 4.003383383383479

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("F7").Value = 25.97000000000062
$ws.Range("H7").Value = 0.00001292125954543977
$ws.Range("I7").Value = 0.00001292125954543977
$ws.Range("L7").Value = 58.95675279055251
$ws.Range("M7").Value = "[28.957953778203432, 88.95555180290158]"
$ws.Range("N7").Value = 0.0002654984500507496
$ws.Range("O7").Value = 0.0002654984500507496
$ws.Range("P7").Value = -0.5912106295043857
$ws.Range("Q7").Value = "[-1.1069475616252324, -0.07547369738353904]"
$ws.Range("R7").Value = 0.02560033529143357
$ws.Range("S7").Value = 0.02560033529143357
$ws.Range("T7").Value = 73.62751606790103
$ws.Range("U7").Value = "[57.88776042809127, 89.36727170771078]"
$ws.Range("V7").Value = 0.000000000003250066882287683
$ws.Range("W7").Value = 0.000000000003250066882287683
$ws.Range("X7").Value = 2.443623623623683
$ws.Range("Y7").Value = 0.3119519519519613
$ws.Range("Z7").Value = 4.575295295295405

# Row 8
$ws.Range("F8").Value = 25.97000000000062
$ws.Range("H8").Value = 0.0001541061788941089
$ws.Range("I8").Value = 0.0001541061788941089
$ws.Range("L8").Value = 50.75187419528145
$ws.Range("M8").Value = "[19.26995262010422, 82.23379577045867]"
$ws.Range("N8").Value = 0.00220674558812517
$ws.Range("O8").Value = 0.00220674558812517
$ws.Range("P8").Value = -0.7170001251436169
$ws.Range("Q8").Value = "[-1.2956318050840787, -0.13836844520315505]"
$ws.Range("R8").Value = 0.01630266950170522
$ws.Range("S8").Value = 0.01630266950170522
$ws.Range("T8").Value = 56.66591671683103
$ws.Range("U8").Value = "[40.7329740739031, 72.59885935975896]"
$ws.Range("V8").Value = 0.00000000583978465584778
$ws.Range("W8").Value = 0.00000000583978465584778
$ws.Range("X8").Value = 2.963543543543615
$ws.Range("Y8").Value = 0.5719119119119287
$ws.Range("Z8").Value = 5.355175175175301

# Row 9
$ws.Range("F9").Value = 25.15000000000049
$ws.Range("H9").Value = 0.00008546725352043261
$ws.Range("I9").Value = 0.00008546725352043261
$ws.Range("L9").Value = 51.4641115454054
$ws.Range("M9").Value = "[22.550296062310835, 80.37792702849997]"
$ws.Range("N9").Value = 0.0008256520003879775
$ws.Range("O9").Value = 0.0008256520003879775
$ws.Range("P9").Value = -1.119526511189155
$ws.Range("Q9").Value = "[-1.7233160902574642, -0.5157369321208467]"
$ws.Range("R9").Value = 0.0005271026078337293
$ws.Range("S9").Value = 0.0005271026078337293
$ws.Range("T9").Value = 67.49945576877728
$ws.Range("U9").Value = "[51.94938006663972, 83.04953147091484]"
$ws.Range("V9").Value = 0.00000000002934985587899064
$ws.Range("W9").Value = 0.00000000002934985587899064
$ws.Range("X9").Value = 4.48118118118127
$ws.Range("Y9").Value = 2.064364364364406
$ws.Range("Z9").Value = 6.897997997998135

# Row 10
$ws.Range("F10").Value = 25.15000000000049
$ws.Range("H10").Value = 0.003659314374704237
$ws.Range("I10").Value = 0.003659314374704237
$ws.Range("L10").Value = 41.98735553859395
$ws.Range("M10").Value = "[13.68992514043498, 70.28478593675293]"
$ws.Range("N10").Value = 0.004529329455255082
$ws.Range("O10").Value = 0.004529329455255082
$ws.Range("P10").Value = -1.547210796362541
$ws.Range("Q10").Value = "[-2.4403162154010807, -0.6541053773240013]"
$ws.Range("R10").Value = 0.001095654488354247
$ws.Range("S10").Value = 0.001095654488354247
$ws.Range("T10").Value = 62.42782407981408
$ws.Range("U10").Value = "[45.28102707188113, 79.57462108774703]"
$ws.Range("V10").Value = 0.000000003277776716714698
$ws.Range("W10").Value = 0.000000003277776716714698
$ws.Range("X10").Value = 6.193093093093214
$ws.Range("Y10").Value = 2.61821821821827
$ws.Range("Z10").Value = 9.767967967968158

# Row 11
$ws.Range("F11").Value = 25.15000000000049
$ws.Range("H11").Value = 0.003235624085636712
$ws.Range("I11").Value = 0.003235624085636712
$ws.Range("L11").Value = 42.06438011866459
$ws.Range("M11").Value = "[12.165476420493903, 71.96328381683527]"
$ws.Range("N11").Value = 0.006865770427760376
$ws.Range("O11").Value = 0.006865770427760376
$ws.Range("P11").Value = -1.408842351159387
$ws.Range("Q11").Value = "[-2.264210921506157, -0.5534737808126167]"
$ws.Range("R11").Value = 0.00180489128076422
$ws.Range("S11").Value = 0.00180489128076422
$ws.Range("T11").Value = 60.12060128338094
$ws.Range("U11").Value = "[43.12080813235467, 77.12039443440722]"
$ws.Range("V11").Value = 0.000000006698104071389821
$ws.Range("W11").Value = 0.000000006698104071389821
$ws.Range("X11").Value = 5.639239239239348
$ws.Range("Y11").Value = 2.215415415415459
$ws.Range("Z11").Value = 9.063063063063238
